$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "38.981.19"
$ws.Cells.Item(2,5).Value = "  -5.11%  "

$ws.Cells.Item(3,4).Value = "2.213.64"
$ws.Cells.Item(3,5).Value = "  -7.71%  "

$ws.Cells.Item(4,5).Value = "  -0.02%  "

$ws.Cells.Item(5,4).Value = "'295.75"
$ws.Cells.Item(5,5).Value = "  -5.96%  "

$ws.Cells.Item(6,4).Value = "'80.00"
$ws.Cells.Item(6,5).Value = "  -10.48%  "

$ws.Cells.Item(7,5).Value = "  -5.66%  "

$ws.Cells.Item(8,5).Value = "  +0.01%  "

$ws.Cells.Item(9,5).Value = "  -8.63%  "

$ws.Cells.Item(10,4).Value = "'0.0771"
$ws.Cells.Item(10,5).Value = "  -8.55%  "

$ws.Cells.Item(11,4).Value = "'27.91"
$ws.Cells.Item(11,5).Value = "  -11.43%  "

$ws.Cells.Item(12,4).Value = "'45.90"
$ws.Cells.Item(12,5).Value = "  -13.07%  "

$ws.Cells.Item(13,5).Value = "  -1.67%  "

$ws.Cells.Item(14,4).Value = "2.558.08"
$ws.Cells.Item(14,5).Value = "  -7.57%  "

$ws.Cells.Item(15,4).Value = "'6.07"
$ws.Cells.Item(15,5).Value = "  -8.73%  "

$ws.Cells.Item(16,5).Value = "  -9.71%  "

$ws.Cells.Item(17,4).Value = "2.229.47"
$ws.Cells.Item(17,5).Value = "  -6.04%  "

$ws.Cells.Item(18,5).Value = "  -8.71%  "

$ws.Cells.Item(19,4).Value = "38.894.35"
$ws.Cells.Item(19,5).Value = "  -5.21%  "

$ws.Cells.Item(20,4).Value = "0.0₃0857"
$ws.Cells.Item(20,5).Value = "  -6.83%  "

$ws.Cells.Item(21,4).Value = "'5.68"
$ws.Cells.Item(21,5).Value = "  -8.63%  "

$ws.Cells.Item(22,4).Value = "'64.45"
$ws.Cells.Item(22,5).Value = "  -7.55%  "

$ws.Cells.Item(23,4).Value = "'9.74"
$ws.Cells.Item(23,5).Value = "  -10.92%  "

$ws.Cells.Item(24,4).Value = "'224.33"
$ws.Cells.Item(24,5).Value = "  -4.46%  "

$ws.Cells.Item(25,5).Value = "  -0.08%  "

$ws.Cells.Item(26,5).Value = "  -11.15%  "

$ws.Cells.Item(27,4).Value = "'1.73"
$ws.Cells.Item(27,5).Value = "  -5.72%  "

$ws.Cells.Item(28,4).Value = "'22.22"
$ws.Cells.Item(28,5).Value = "  -8.00%  "

$ws.Cells.Item(29,5).Value = "  -2.82%  "

$ws.Cells.Item(30,4).Value = "'8.87"
$ws.Cells.Item(30,5).Value = "  -6.07%  "

$ws.Cells.Item(31,4).Value = "'148.47"
$ws.Cells.Item(31,5).Value = "  -3.90%  "

$ws.Cells.Item(32,4).Value = "'30.98"
$ws.Cells.Item(32,5).Value = "  -10.28%  "

$ws.Cells.Item(33,5).Value = "  -0.29%  "

$ws.Cells.Item(34,5).Value = "  -10.29%  "

$ws.Cells.Item(35,5).Value = "  -4.97%  "

$ws.Cells.Item(36,5).Value = "  -7.67%  "

$ws.Cells.Item(37,5).Value = "  -5.37%  "

$ws.Cells.Item(38,2).Value = "Kaspa"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(38,4).Value = "'0.0950"
$ws.Cells.Item(38,5).Value = "  -5.84%  "

$ws.Cells.Item(39,2).Value = "LidoDAOToken"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(39,4).Value = "'2.62"
$ws.Cells.Item(39,5).Value = "  -7.66%  "

$ws.Cells.Item(40,2).Value = "ARBITRUM"
$ws.Cells.Item(40,3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(40,4).Value = "'1.58"
$ws.Cells.Item(40,5).Value = "  -9.90%  "

$ws.Cells.Item(41,2).Value = "Celestia"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(41,4).Value = "'14.22"
$ws.Cells.Item(41,5).Value = "  -12.45%  "

$ws.Cells.Item(42,5).Value = "  -7.61%  "

$ws.Cells.Item(43,4).Value = "1.896.85"
$ws.Cells.Item(43,5).Value = "  -4.16%  "

$ws.Cells.Item(44,5).Value = "  -11.58%  "

$ws.Cells.Item(45,5).Value = "  -7.59%  "

$ws.Cells.Item(46,4).Value = "'16.05"
$ws.Cells.Item(46,5).Value = "  -10.34%  "

$ws.Cells.Item(47,4).Value = "'8.81"
$ws.Cells.Item(47,5).Value = "  -9.08%  "

$ws.Cells.Item(48,4).Value = "'2.49"
$ws.Cells.Item(48,5).Value = "  -11.08%  "

$ws.Cells.Item(49,4).Value = "2.426.08"
$ws.Cells.Item(49,5).Value = "  -7.63%  "

$ws.Cells.Item(50,4).Value = "'69.62"
$ws.Cells.Item(50,5).Value = "  -5.23%  "

$ws.Cells.Item(51,2).Value = "TrustWalletToken"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(51,4).Value = "'1.07"
$ws.Cells.Item(51,5).Value = "  -1.55%  "
